# Updated symbol list on Wed Dec 28 17:18:13 UTC 2022 with GitHub Actions
# Applies the cryptocurrency price/volume updates described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a numeric-looking text value while preserving it as
# plain text (so things like trailing zeros / leading zeros survive) and
# without leaving a lingering custom number-format style on the cell.
function Set-TextValue($cell, [string]$text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Column D (Price) updates
Set-TextValue $ws.Range("D2")  "243.35"
Set-TextValue $ws.Range("D3")  "23.78"
Set-TextValue $ws.Range("D4")  "5.261"
Set-TextValue $ws.Range("D5")  "0.05819"
Set-TextValue $ws.Range("D6")  "6.464"
Set-TextValue $ws.Range("D9")  "0.9018"
Set-TextValue $ws.Range("D10") "0.1379"
Set-TextValue $ws.Range("D11") "0.07083"
Set-TextValue $ws.Range("D12") "0.03086"
Set-TextValue $ws.Range("D13") "0.03033"
Set-TextValue $ws.Range("D14") "0.09318"
Set-TextValue $ws.Range("D15") "3.819"
Set-TextValue $ws.Range("D16") "0.001562"
Set-TextValue $ws.Range("D17") "0.04698"
Set-TextValue $ws.Range("D18") "0.0006020"
Set-TextValue $ws.Range("D19") "0.006241"
Set-TextValue $ws.Range("D20") "0.001260"
Set-TextValue $ws.Range("D21") "0.003879"
Set-TextValue $ws.Range("D22") "0.00008708"
Set-TextValue $ws.Range("D23") "3.558"
Set-TextValue $ws.Range("D24") "2.173"
Set-TextValue $ws.Range("D25") "0.3192"
Set-TextValue $ws.Range("D26") "0.1317"
Set-TextValue $ws.Range("D28") "0.0002328"
Set-TextValue $ws.Range("D40") "0.03787"
Set-TextValue $ws.Range("D41") "0.006306"
Set-TextValue $ws.Range("D42") "0.1051"
Set-TextValue $ws.Range("D43") "0.002524"
Set-TextValue $ws.Range("D44") "0.006941"
Set-TextValue $ws.Range("D45") "0.00005311"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("D48") "0.007729"
Set-TextValue $ws.Range("D49") "0.00002102"
Set-TextValue $ws.Range("D50") "0.0002002"

# Column E (Volume(1h)) label updates
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
